$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep column D as text so numeric-looking price strings (e.g. "27.455.36",
# "1.320") are not auto-converted into numbers by Excel.
$ws.Columns.Item(4).NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
$ws.Range("D2").Value = "27.455.36"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.825.17"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.82%  "
$ws.Range("D5").Value = "331.57"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D7").Value = "0.4542"
$ws.Range("E7").Value = "  -1.92%  "
$ws.Range("D8").Value = "0.3804"
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("D9").Value = "46.26"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").Value = "0.07883"
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("D11").Value = "0.9702"
$ws.Range("E11").Value = "  -2.43%  "
$ws.Range("D12").Value = "20.99"
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("D13").Value = "1.830.58"
$ws.Range("E13").Value = "  -2.12%  "
$ws.Range("D14").Value = "5.859"
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").Value = "7.037"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("D16").Value = "1.007"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").Value = "88.53"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "0.06632"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").Value = "0.00001026"
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("D20").Value = "17.14"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("D21").Value = "1.008"
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("D22").Value = "27.435.36"
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("D23").Value = "5.319"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("D24").Value = "10.78"
$ws.Range("E24").Value = "  -0.60%  "
$ws.Range("D25").Value = "2.304"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("D26").Value = "2.054.88"
$ws.Range("E26").Value = "  -1.65%  "
$ws.Range("D27").Value = "156.58"
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("D28").Value = "19.38"
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("D29").Value = "2.060"
$ws.Range("E29").Value = "  -2.36%  "
$ws.Range("D30").Value = "5.234"
$ws.Range("E30").Value = "  -2.02%  "
$ws.Range("D32").Value = "0.9438"
$ws.Range("E32").Value = "  -2.39%  "
$ws.Range("D33").Value = "0.09286"
$ws.Range("E33").Value = "  -1.43%  "
$ws.Range("D34").Value = "3.575"
$ws.Range("E34").Value = "  -1.85%  "
$ws.Range("D35").Value = "5.230"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").Value = "1.320"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "0.05916"
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("D38").Value = "0.02179"
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("D39").Value = "1.160"
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("D40").Value = "8.017"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("D41").Value = "0.5739"
$ws.Range("E41").Value = "  -2.42%  "
$ws.Range("D42").Value = "0.1830"
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("D43").Value = "9.993"
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("D44").Value = "1.259"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").Value = "0.5449"
$ws.Range("E45").Value = "  -2.60%  "
$ws.Range("D46").Value = "11.94"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").Value = "1.862"
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("D48").Value = "0.06613"
$ws.Range("E48").Value = "  -2.19%  "
$ws.Range("D49").Value = "110.19"
$ws.Range("E49").Value = "  -1.54%  "
$ws.Range("D50").Value = "1.037"
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("D51").Value = "1.005"
$ws.Range("E51").Value = "  -0.67%  "

# Rows 45 and 46 swapped coin identities (Decentraland <-> EnergySwap),
# including their links, while keeping row-rank (column A) unchanged.
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
